# PAS-6576 - Update "individual VIN retrieval" logic to use ENTRY DATE and VALID
# Reworks the VIN upload test fixture: row 3 becomes an "invalid VIN" sample
# (ENTRYDATE 2001-01-01, VALID=N), row 4 is filled in as a second full sample
# row, and a brand new row 5 ("secondValid") is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 3: mark this VIN row as the "invalid" sample and change its entry date
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = "invalidVin"
$ws.Range("F3").Value = "invalidVin"
$ws.Range("AI3").Value = 20010101
$ws.Range("AJ3").Value = "N"

# ---------------------------------------------------------------------------
# Row 4: bring over row 3's formatting, then populate every column with the
# second full data row (previously this row only had empty, styled cells).
# ---------------------------------------------------------------------------
$ws.Range("A3:AL3").Copy()
$ws.Range("A4:AL4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A4").Value = "1HGEM215&4"
$ws.Range("B4").Value = "SYMBOL_2000_SS_TEST"
$ws.Range("C4").Value = 2005
$ws.Range("D4").Value = "TEST"
$ws.Range("E4").Value = "TEST"
$ws.Range("F4").Value = "TEST"
$ws.Range("G4").Value = "MDX ADVANCE"
$ws.Range("H4").Value = 53080
$ws.Range("I4").Value = "WAG"
$ws.Range("J4").Value = "TEST"
$ws.Range("K4").Value = "TEST"
$ws.Range("L4").Value = "TEST"
$ws.Range("M4").Value = "WAG"
$ws.Range("N4").Value = "3.5L V6   "
$ws.Range("O4").Value = 6
$ws.Range("P4").Value = "G"
$ws.Range("Q4").Value = 214
$ws.Range("R4").Value = "2WD"
$ws.Range("S4").Value = 2
$ws.Range("T4").Value = "000R"
$ws.Range("U4").Value = "DUAL AIR BAGS FRONT"
$ws.Range("V4").Value = 2
$ws.Range("W4").Value = "4 WHEEL STANDARD"
$ws.Range("X4").Value = "STD"
$ws.Range("Y4").Value = "B-IMMOBILIZER/KEYLSS ENTRY/ALARM"
$ws.Range("Z4").Value = "I"
$ws.Range("AA4").Value = 39
$ws.Range("AB4").Value = 40
$ws.Range("AC4").Value = "K"
$ws.Range("AD4").Value = "Y"
$ws.Range("AE4").Value = "X"
$ws.Range("AF4").Value = "X"
$ws.Range("AG4").Value = "X"
$ws.Range("AH4").Value = "X"
$ws.Range("AI4").Value = 20020101
$ws.Range("AJ4").Value = "Y"
$ws.Range("AK4").Value = "Y"
$ws.Range("AL4").Value = "N"

# ---------------------------------------------------------------------------
# Row 5 (new): a second "valid" VIN row, same shape as row 4 but with its own
# MAKE_TEXT/MODEL_TEXT ("secondValid") and entry date.
# ---------------------------------------------------------------------------
$ws.Range("A3:AL3").Copy()
$ws.Range("A5:AL5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A5").Value = "1HGEM215&4"
$ws.Range("B5").Value = "SYMBOL_2000_SS_TEST"
$ws.Range("C5").Value = 2005
$ws.Range("D5").Value = "TEST"
$ws.Range("E5").Value = "secondValid"
$ws.Range("F5").Value = "secondValid"
$ws.Range("G5").Value = "MDX ADVANCE"
$ws.Range("H5").Value = 53080
$ws.Range("I5").Value = "WAG"
$ws.Range("J5").Value = "TEST"
$ws.Range("K5").Value = "TEST"
$ws.Range("L5").Value = "TEST"
$ws.Range("M5").Value = "WAG"
$ws.Range("N5").Value = "3.5L V6   "
$ws.Range("O5").Value = 6
$ws.Range("P5").Value = "G"
$ws.Range("Q5").Value = 214
$ws.Range("R5").Value = "2WD"
$ws.Range("S5").Value = 2
$ws.Range("T5").Value = "000R"
$ws.Range("U5").Value = "DUAL AIR BAGS FRONT"
$ws.Range("V5").Value = 2
$ws.Range("W5").Value = "4 WHEEL STANDARD"
$ws.Range("X5").Value = "STD"
$ws.Range("Y5").Value = "B-IMMOBILIZER/KEYLSS ENTRY/ALARM"
$ws.Range("Z5").Value = "I"
$ws.Range("AA5").Value = 39
$ws.Range("AB5").Value = 40
$ws.Range("AC5").Value = "K"
$ws.Range("AD5").Value = "Y"
$ws.Range("AE5").Value = "X"
$ws.Range("AF5").Value = "X"
$ws.Range("AG5").Value = "X"
$ws.Range("AH5").Value = "X"
$ws.Range("AI5").Value = 20030101
$ws.Range("AJ5").Value = "Y"
$ws.Range("AK5").Value = "Y"
$ws.Range("AL5").Value = "N"

# ---------------------------------------------------------------------------
# Last touched cell in the sheet, as left by the author before saving.
# ---------------------------------------------------------------------------
$ws.Range("J12").Select()
